# Disposal.xlsx edit: change reporting period, update the disposal line's
# amounts/texts, and append a second disposal line for the new period.
#
# Notes on quirks of this sheet / engine that shape the approach below:
#
# * Every data cell (row 7/8) is stored as a shared string, even the
#   numeric-looking ones ("30004", "249000", ...), not as a native number.
#   A plain `Range.Value = "70005"` assignment auto-converts to a number
#   (General number format), which would also fork a brand-new cell style.
#   To keep a numeric-looking value as text -- and keep using the existing
#   cell style -- the value is written into a scratch cell via a formula
#   (`="70005"`), copied, then Paste-Special "Values only" (xlPasteValues)
#   into the destination: a formula's result pastes as a plain text/number
#   *value*, landing a shared-string cell without touching formatting.
#
# * That Paste-Special trick unmerges any merged destination cell as a
#   side effect, and re-merging afterwards forks the border styles of the
#   (now touched) follower cells. Row 2 is merged (A2:G2) but its new text
#   isn't numeric-looking, so it's set with a plain `.Value =`, which
#   leaves merges/styles untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("ZZ1")

function Set-TextValue {
    param($range, [string]$text)

    $escaped = $text.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

# --- Reporting period (row 2, merged A2:G2) --------------------------------
$ws.Range("A2").Value = "Период: 2023-11-01 - 2023-11-30"

# --- Update the existing disposal row (row 8) -------------------------------
Set-TextValue $ws.Range("A8") "70005"
$ws.Range("B8").Value = "Сбербанк"
# C8 ("лебенков") and F8 ("Выбытие") are unchanged.
Set-TextValue $ws.Range("D8") "5000"
$ws.Range("E8").Value = "01.11.2023 00:00:00"
$ws.Range("G8").Value = "Докторская"

# --- Append the new disposal row (row 9), same data as row 8 except for
# the order number in column A -----------------------------------------------
$ws.Range("B9").Value = "Сбербанк"
$ws.Range("C9").Value = "лебенков"
$ws.Range("E9").Value = "01.11.2023 00:00:00"
$ws.Range("F9").Value = "Выбытие"
$ws.Range("G9").Value = "Докторская"
$ws.Range("A9:G9").Borders.LineStyle = 1
Set-TextValue $ws.Range("D9") "5000"
Set-TextValue $ws.Range("A9") "70006"

$scratch.ClearContents()
$excel.CutCopyMode = $false

# --- Column G got a bit wider to fit the new text ---------------------------
$ws.Columns.Item(7).ColumnWidth = 12.0177928379604
